$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 510, pushing all existing
# data (rows 510-619) down to rows 512-621. This mirrors the weekly
# refresh pattern used throughout this sheet: the newest week's prices
# are always prepended right after the most recent previously-recorded
# date block.
$ws.Rows.Item(510).Insert()
$ws.Rows.Item(510).Insert()

# New row 510: "Primera" quality row for the new date (2022-03-21 / 44641)
$ws.Cells.Item(510, 1).Value = 8
$ws.Cells.Item(510, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(510, 3).Value = "Coquimbo"
$ws.Cells.Item(510, 4).Value = 44641
$ws.Cells.Item(510, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(510, 5).Value = 4
$ws.Cells.Item(510, 6).Value = 100112023
$ws.Cells.Item(510, 7).Value = "Brócoli"
$ws.Cells.Item(510, 8).Value = "Sin especificar"
$ws.Cells.Item(510, 9).Value = "Primera"
$ws.Cells.Item(510, 10).Value = 2500
$ws.Cells.Item(510, 11).Value = 950
$ws.Cells.Item(510, 12).Value = 1000
$ws.Cells.Item(510, 13).Value = 975
$ws.Cells.Item(510, 14).Value = "`$/unidad"
$ws.Cells.Item(510, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(510, 16).Value = 975
$ws.Cells.Item(510, 17).Value = 1
$ws.Cells.Item(510, 18).Value = "Hortaliza"

# New row 511: "Segunda" quality row for the same new date
$ws.Cells.Item(511, 1).Value = 8
$ws.Cells.Item(511, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(511, 3).Value = "Coquimbo"
$ws.Cells.Item(511, 4).Value = 44641
$ws.Cells.Item(511, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(511, 5).Value = 4
$ws.Cells.Item(511, 6).Value = 100112023
$ws.Cells.Item(511, 7).Value = "Brócoli"
$ws.Cells.Item(511, 8).Value = "Sin especificar"
$ws.Cells.Item(511, 9).Value = "Segunda"
$ws.Cells.Item(511, 10).Value = 1240
$ws.Cells.Item(511, 11).Value = 850
$ws.Cells.Item(511, 12).Value = 900
$ws.Cells.Item(511, 13).Value = 875
$ws.Cells.Item(511, 14).Value = "`$/unidad"
$ws.Cells.Item(511, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(511, 16).Value = 875
$ws.Cells.Item(511, 17).Value = 1
$ws.Cells.Item(511, 18).Value = "Hortaliza"
